# 2.a.1.xlsx — add the "2020" data point as a new column Q next to the
# existing 2007-2019 series on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (year headers): Q3 = 2020, formatted like the preceding P3 (2019) cell
$ws.Range("P3").Copy() | Out-Null
$ws.Range("Q3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("Q3").Value = 2020

# --- Row 4 (data values): Q4 = new 2020 index value, formatted like P4 (2019)
$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("Q4").Value = 0.067156049127444606

# --- Tighten the data row's number format from the custom "0.0" format to
#     the built-in two-decimal "0.00" format, across the whole series
#     (including the newly added Q4 cell).
$ws.Range("D4:Q4").NumberFormat = "0.00"

# --- Clear the stale selection left over from a previous editing session
#     (was parked on B12, well outside the used range) back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
